$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values contain dotted numeric-looking text (e.g. "67.634.79" or "595.29").
# Force text format so Excel does not reinterpret them as numbers, then restore the default
# "Normal" style so no stray number formatting is left behind on the cells.
$dCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D25", "D26", "D31", "D32", "D34", "D36", "D38", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.634.79'
$ws.Range("D3").Value = '3.781.76'
$ws.Range("D5").Value = '595.29'
$ws.Range("D6").Value = '166.42'
$ws.Range("D7").Value = '3.781.00'
$ws.Range("D9").Value = '0.520'
$ws.Range("D14").Value = '36.27'
$ws.Range("D15").Value = '4.417.26'
$ws.Range("D16").Value = '3.774.44'
$ws.Range("D17").Value = '18.45'
$ws.Range("D18").Value = '67.611.08'
$ws.Range("D20").Value = '6.99'
$ws.Range("D21").Value = '10.08'
$ws.Range("D22").Value = '457.41'
$ws.Range("D25").Value = '83.39'
$ws.Range("D26").Value = '11.94'
$ws.Range("D31").Value = '7.28'
$ws.Range("D32").Value = '29.82'
$ws.Range("D34").Value = '9.21'
$ws.Range("D36").Value = '3.735.40'
$ws.Range("D38").Value = '3.35'
$ws.Range("D40").Value = '0.992'
$ws.Range("D41").Value = '5.76'
$ws.Range("D42").Value = '1.00'
$ws.Range("D44").Value = '45.36'
$ws.Range("D46").Value = '47.12'
$ws.Range("D47").Value = '8.34'
$ws.Range("D48").Value = '148.34'
$ws.Range("D50").Value = '389.99'
$ws.Range("D51").Value = '25.56'

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Column E (Volume 1h) values are plain text percentages padded with spaces; a direct
# assignment keeps them as text since Excel will not parse the padded string as a number.
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  -2.03%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E17").Value = '  +3.40%  '
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").Value = '  -5.71%  '
$ws.Range("E22").Value = '  -2.15%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("E24").Value = '  +8.12%  '
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  +5.88%  '
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("E46").Value = '  +2.79%  '
$ws.Range("E47").Value = '  -2.87%  '
$ws.Range("E48").Value = '  +1.11%  '
$ws.Range("E49").Value = '  -4.16%  '
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("E51").Value = '  +0.70%  '
